{"js": "// Center-align the paragraphs inside every cell of the document's\n// (first) table \u2014 the \"Kode\" / \"Penjelasan\" header row plus the blank\n// data row beneath it \u2014 matching the author's table-formatting fix.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (const table of tables.items) {\n  const rows = table.rows;\n  rows.load(\"items\");\n  await context.sync();\n\n  for (const row of rows.items) {\n    const cells = row.cells;\n    cells.load(\"items\");\n    await context.sync();\n\n    for (const cell of cells.items) {\n      const paragraphs = cell.body.paragraphs;\n      paragraphs.load(\"items\");\n      await context.sync();\n\n      for (const paragraph of paragraphs.items) {\n        paragraph.alignment = Word.Alignment.centered;\n      }\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Center-align the paragraphs inside every cell of the document's table\n# (the \"Kode\" / \"Penjelasan\" header row plus the blank data row beneath\n# it) to match the author's table-formatting fix.\n$d = $word.ActiveDocument\n\n# wdAlignParagraphCenter = 1\n$wdAlignParagraphCenter = 1\n\nforeach ($tbl in $d.Tables) {\n    foreach ($cell in $tbl.Range.Cells) {\n        foreach ($para in $cell.Range.Paragraphs) {\n            $para.Alignment = $wdAlignParagraphCenter\n        }\n    }\n}\n"}
